$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "Mengyao Piao"
$ws.Range("B5").Value = "Chaoda Song"
$ws.Range("B6").Value = "Jun Liu"

$ws.Range("B10").Value = "https://github.com/upc2017/assignments-"
$ws.Range("D10").Value = "Yes"

$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 4
$ws.Range("D15").Value = 4

$ws.Range("D20").Value = 4
$ws.Range("D21").Value = 4
$ws.Range("D22").Value = 4
$ws.Range("D23").Value = 4
$ws.Range("D24").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("D27").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("D29").Value = 4
$ws.Range("D30").Value = 4

$ws.Range("D33").Value = 4
$ws.Range("D34").Value = 4
$ws.Range("D35").Value = 4
$ws.Range("D36").Value = 4
$ws.Range("D37").Value = 4
$ws.Range("D38").Value = 4
$ws.Range("D39").Value = 4
$ws.Range("D40").Value = 4
$ws.Range("D41").Value = 4

$ws.Range("D46").Value = 4
$ws.Range("D47").Value = 4
$ws.Range("D48").Value = 4
$ws.Range("D49").Value = 4
$ws.Range("D50").Value = 4

$ws.Range("D50").Select()
